$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "571.41") must be
# forced to remain TEXT (matching the source inlineStr cells) instead of
# being auto-coerced to a numeric type by Excel's type inference.
function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.086.51"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.509.01"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "571.41"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextValue "D6" "184.45"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("D8").Value = "3.503.86"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("E11").Value = "  -2.15%  "
Set-TextValue "D12" "54.10"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("E13").Value = "  +0.04%  "
Set-TextValue "D14" "9.43"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "4.073.58"
$ws.Range("E15").Value = "  -1.77%  "
Set-TextValue "D16" "19.34"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "69.009.13"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "3.501.52"
$ws.Range("E18").Value = "  -2.02%  "
Set-TextValue "D19" "12.27"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -1.28%  "
Set-TextValue "D21" "539.80"
$ws.Range("E21").Value = "  +13.75%  "
$ws.Range("E22").Value = "  -2.88%  "
Set-TextValue "D23" "18.74"
$ws.Range("E23").Value = "  -2.49%  "
Set-TextValue "D24" "5.00"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +0.61%  "
Set-TextValue "D26" "93.92"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  -2.96%  "
Set-TextValue "D28" "10.81"
$ws.Range("E28").Value = "  -1.93%  "
Set-TextValue "D29" "9.15"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("E30").Value = "  -1.83%  "
Set-TextValue "D31" "7.26"
$ws.Range("E31").Value = "  -9.24%  "
Set-TextValue "D32" "12.57"
$ws.Range("E32").Value = "  +2.71%  "
Set-TextValue "D33" "64.78"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -4.17%  "
Set-TextValue "D35" "565.47"
$ws.Range("E35").Value = "  -2.73%  "
Set-TextValue "D36" "37.91"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").Value = "0.0₃0767"
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -2.87%  "
Set-TextValue "D43" "3.35"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D44" "2.99"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "3.228.62"
$ws.Range("E46").Value = "  -0.18%  "
Set-TextValue "D47" "0.0441"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -4.81%  "
Set-TextValue "D50" "0.997"
$ws.Range("E50").Value = "  -0.22%  "
Set-TextValue "D51" "138.00"
$ws.Range("E51").Value = "  +2.60%  "
